$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from E1 (header style) into F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:52:08.813026"
$ws.Range("F3").Value = "2021-10-05 10:52:08.813037"
$ws.Range("F4").Value = "2021-10-05 10:52:08.813041"
$ws.Range("F5").Value = "2021-10-05 10:52:08.813044"
$ws.Range("F6").Value = "2021-10-05 10:52:08.813047"
$ws.Range("F7").Value = "2021-10-05 10:52:08.813050"
$ws.Range("F8").Value = "2021-10-05 10:52:08.813052"
$ws.Range("F9").Value = "2021-10-05 10:52:08.813055"
$ws.Range("F10").Value = "2021-10-05 10:52:08.813058"
$ws.Range("F11").Value = "2021-10-05 10:52:08.813061"
$ws.Range("F12").Value = "2021-10-05 10:52:08.813064"
$ws.Range("F13").Value = "2021-10-05 10:52:08.813067"
$ws.Range("F14").Value = "2021-10-05 10:52:08.813069"
$ws.Range("F15").Value = "2021-10-05 10:52:08.813072"
$ws.Range("F16").Value = "2021-10-05 10:52:08.813075"
$ws.Range("F17").Value = "2021-10-05 10:52:08.813078"
$ws.Range("F18").Value = "2021-10-05 10:52:08.813081"
$ws.Range("F19").Value = "2021-10-05 10:52:08.813084"
$ws.Range("F20").Value = "2021-10-05 10:52:08.813087"
$ws.Range("F21").Value = "2021-10-05 10:52:08.813089"
$ws.Range("F22").Value = "2021-10-05 10:52:08.813092"
$ws.Range("F23").Value = "2021-10-05 10:52:08.813095"
$ws.Range("F24").Value = "2021-10-05 10:52:08.813098"
$ws.Range("F25").Value = "2021-10-05 10:52:08.813100"
$ws.Range("F26").Value = "2021-10-05 10:52:08.813103"
$ws.Range("F27").Value = "2021-10-05 10:52:08.813106"
$ws.Range("F28").Value = "2021-10-05 10:52:08.813109"
$ws.Range("F29").Value = "2021-10-05 10:52:08.813112"
$ws.Range("F30").Value = "2021-10-05 10:52:08.813114"
$ws.Range("F31").Value = "2021-10-05 10:52:08.813117"
$ws.Range("F32").Value = "2021-10-05 10:52:08.813120"
$ws.Range("F33").Value = "2021-10-05 10:52:08.813123"
$ws.Range("F34").Value = "2021-10-05 10:52:08.813126"
$ws.Range("F35").Value = "2021-10-05 10:52:08.813129"
$ws.Range("F36").Value = "2021-10-05 10:52:08.813131"
$ws.Range("F37").Value = "2021-10-05 10:52:08.813134"
$ws.Range("F38").Value = "2021-10-05 10:52:08.813137"
$ws.Range("F39").Value = "2021-10-05 10:52:08.813140"
$ws.Range("F40").Value = "2021-10-05 10:52:08.813143"
$ws.Range("F41").Value = "2021-10-05 10:52:08.813146"
$ws.Range("F42").Value = "2021-10-05 10:52:08.813149"
